$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.32297956943512
$ws.Range("B1").Value = 1.550219416618347
$ws.Range("C1").Value = 1.990261435508728
$ws.Range("D1").Value = 2.716326236724854
$ws.Range("E1").Value = 6.544262409210205
